$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Schedule_date")

# --- Fix "2: Coding fundamental" -> "2: Coding fundamentals" typo on the Schedule sheet ---
$ws1.Range("B4").Value = "2: Coding fundamentals"
$ws1.Range("B5").Value = "2: Coding fundamentals"
$ws1.Range("B6").Value = "2: Coding fundamentals"

# --- Update topics: "Grammar of graphics" -> "ggplot 101", "Themes, labels, facets" -> "Themes, labels, facets (ggplot 102)" ---
$ws1.Range("C5").Value = "ggplot 101"
$ws1.Range("C6").Value = "Themes, labels, facets (ggplot 102)"

$ws2.Range("D5").Value = "ggplot 101"
$ws2.Range("D6").Value = "Themes, labels, facets (ggplot 102)"

# --- Widen column C on Schedule_date (drop AutoFit/bestFit, set an explicit width) ---
$ws2.Columns.Item(3).ColumnWidth = 21.6666666666667

# --- Update selections / active sheet to match new cursor positions ---
[void]$ws2.Range("D5:D6").Select()
[void]$ws1.Activate()
[void]$ws1.Range("C5:C6").Select()

$wb.Save()
